$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial value for every data row (2-461).
# The commit updates that "changed on" date from 45172 (2023-09-03) to
# 45175 (2023-09-06) for every row.
$lastRow = 461
$ws.Range("C2:C$lastRow").Value = 45175
